$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1. Remove the "text" category column (column Y) entirely. This shifts all
#    columns to its right (Z, AA, AB, AC, AD, AE) one position to the left
#    (Z->Y, AA->Z, AB->AA, AC->AB, AD->AC, AE->AD).
# ---------------------------------------------------------------------------
$ws.Columns("Y").Delete()

# ---------------------------------------------------------------------------
# 2. Add the new "storeKeys(json,jsonpath,var)" expression to the `json`
#    function list (column M), inserted alphabetically just before
#    "storeValue(json,jsonpath,var)". Since the json list only occupied
#    M2:M17 previously (now needs to become M2:M18), we just rewrite the
#    tail of the list directly.
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 13).Value = "storeKeys(json,jsonpath,var)"
$ws.Cells.Item(17, 13).Value = "storeValue(json,jsonpath,var)"
$ws.Cells.Item(18, 13).Value = "storeValues(json,jsonpath,var)"

# ---------------------------------------------------------------------------
# 3. The "text" category itself is also removed from the `target` list
#    (column A). That list has "text" at A25; remove it and shift the
#    remaining entries (web, webalert, webcookie, ws, ws.async, xml) up by
#    one row, clearing the now-unused last row (A31).
# ---------------------------------------------------------------------------
for ($r = 26; $r -le 31; $r++) {
    $v = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r - 1, 1).Value = $v
}
$ws.Cells.Item(31, 1).ClearContents()

# ---------------------------------------------------------------------------
# 4. Update the named ranges affected by the edits above.
# ---------------------------------------------------------------------------
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
